$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update F1: last status check time 02:30 -> 02:45
$ws.Range("F1").Value = "Last status check on: 06.02.2022 02:45"

# Update D3: "+0.4" text -> numeric 0.4
$ws.Range("D3").Value = 0.4

# Update E3: inline string date -> numeric Excel date serial with datetime formatting
$ws.Range("E3").Value = 44598.10496527778
$ws.Range("E3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
